$d = $word.ActiveDocument

# 1) Update both "Date de création" / "Date de version" table cells:
#    10/07/2018 -> 11/07/2018 (both occurrences in the document change identically)
$d.Content.Find.Execute("10/07/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "11/07/2018", 2)

# 2) Swap the order of the two bullet items in the "eCollection" row
#    (numId 47): "Ajout de la fonction filtrage..." and
#    "Mise en place du bridge entre Formbuilder et Ecollection" trade places.
#    There is a second, unrelated occurrence of the bridge sentence further
#    down (FormBuilder row, numId 49) which must stay untouched, so locate
#    the specific pair of paragraphs by their current text instead of doing
#    a blanket find/replace.
$textFiltrage = "Ajout de la fonction filtrage et maquettage des diverses pop-up modal"
$textBridge = "Mise en place du bridge entre Formbuilder et Ecollection"

$pFiltrage = $null
$pBridge = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $textFiltrage) {
        $pFiltrage = $p
    }
}

if ($pFiltrage -ne $null) {
    # The bridge paragraph that needs swapping is the one immediately
    # following the filtrage paragraph.
    $pBridge = $pFiltrage.Next()
}

if (($pFiltrage -ne $null) -and ($pBridge -ne $null)) {
    $pFiltrage.Range.Text = $textBridge
    $pBridge.Range.Text = $textFiltrage
}
